$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Session 1 (2023-12-11 09:00:05) - fix wording/punctuation on three items
$ws.Range("B7").Value  = "Ich würde mich unbehaglich dabei fühlen, wenn ich einen Job bekommen würde, bei dem ich mit Robotern arbeiten müsste."
$ws.Range("B5").Value  = "Ich würde mich im Gespräch mit einem Roboter entspannt fühlen."
$ws.Range("B14").Value = "Die Zuweisung von Routineaufgaben an Roboter ermöglicht es den Menschen, bedeutungsvollere Aufgaben zu erledigen."

# Session 2 (2023-12-11 09:01:48) - trim trailing spaces on two items
$ws.Range("B6").Value = "Wenn Roboter Emotionen hätten, könnte ich mich mit ihnen anfreunden."
$ws.Range("B8").Value = "Ich befürchte, dass Roboter meine Anweisungen nicht verstehen wüden."

# Session 3 (2023-12-11 09:02:53) - trim trailing spaces on three more items
$ws.Range("B11").Value = "Ich möchte nicht, dass ein Roboter mich anfasst."
$ws.Range("B12").Value = "Roboter sind notwendig, da sie Jobs ausführen können, welche zu schwer oder gefährlich für Menschen sind."
$ws.Range("B13").Value = "Roboter können das Leben einfacher machen."

# Final cursor/selection position left by the author
$ws.Range("F18").Select()
